$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 87 (existing rows 87.. shift down to 90..)
$ws.Rows.Item(87).Insert()
$ws.Rows.Item(87).Insert()
$ws.Rows.Item(87).Insert()

# Common values shared by these Espárragos / Lo Valledor rows
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$catId = 300000000
$categoria = "Espárragos"
$unidad = "`$/kilo"
$kgUnidades = 1
$clasificacion = "Hortaliza"

function Set-FilaEsparrago($Fila, $Fecha, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Origen) {
    $ws.Cells.Item($Fila, 1).Value = $mercadoId
    $ws.Cells.Item($Fila, 2).Value = $mercado
    $ws.Cells.Item($Fila, 3).Value = $region
    $ws.Cells.Item($Fila, 4).Value = $Fecha
    $ws.Cells.Item($Fila, 5).Value = $codreg
    $ws.Cells.Item($Fila, 6).Value = $catId
    $ws.Cells.Item($Fila, 7).Value = $categoria
    $ws.Cells.Item($Fila, 8).Value = $Variedad
    $ws.Cells.Item($Fila, 9).Value = $Calidad
    $ws.Cells.Item($Fila, 10).Value = $Volumen
    $ws.Cells.Item($Fila, 11).Value = $PrecioMin
    $ws.Cells.Item($Fila, 12).Value = $PrecioMax
    $ws.Cells.Item($Fila, 13).Value = $PrecioProm
    $ws.Cells.Item($Fila, 14).Value = $unidad
    $ws.Cells.Item($Fila, 15).Value = $Origen
    $ws.Cells.Item($Fila, 16).Value = $PrecioProm
    $ws.Cells.Item($Fila, 17).Value = $kgUnidades
    $ws.Cells.Item($Fila, 18).Value = $clasificacion
}

Set-FilaEsparrago 87 44518 "Sin especificar" "Banquete" 550 1300 1400 1358 "Provincia de Linares"
Set-FilaEsparrago 88 44518 "Sin especificar" "Primera" 460 1100 1200 1157 "Provincia de Linares"
Set-FilaEsparrago 89 44518 "Sin especificar" "Segunda" 300 900 1000 950 "Provincia de Linares"
